# Apply edits described in the commit:
# "Se han añadido nuevos inputs y modificación restricción soft 2"
#
# This updates several manually-entered input cells on Sheet1:
#  - "1. Horas y profesor para cada clase y asignatura" block (rows 24-27)
#  - "2. Disponibilidad profesor" block (rows 74-82)
#  - the "Horas totales" row (row 90)
# All other changed cells in the original diff are formulas that will
# recalculate automatically once these inputs change.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Horas totales (row 90): Friday (H) 13 -> 14 ---
# (Set this before the availability cells below so every downstream
# formula that depends on H90 picks up the new total right away.)
$ws.Range("H90").Value = 14

# --- Block 1: Horas y profesor para cada clase y asignatura (rows 24-27) ---
# Clase A (row 24): Tutoria hours 1 -> 0, Quimica hours 0 -> 1
$ws.Range("H24").Value = 0
$ws.Range("L24").Value = 1

# Clase B (row 25): Tutoria hours 1 -> 0, Ingles hours 3 -> 4
$ws.Range("H25").Value = 0
$ws.Range("X25").Value = 4

# Clase C (row 26): Tutoria hours 1 -> 0, Filosofia hours 0 -> 1
$ws.Range("H26").Value = 0
$ws.Range("P26").Value = 1

# Clase D (row 27): Matematicas hours 4 -> 5, Tutoria hours 1 -> 0
$ws.Range("D27").Value = 5
$ws.Range("H27").Value = 0

# --- Block 2: Disponibilidad profesor (rows 74-82) ---
# MARIO (row 74): clear L (unavailable) flag
$ws.Range("D74").Value = ""

# MARIA (row 75): clear M and V flags
$ws.Range("E75").Value = ""
$ws.Range("H75").Value = ""

# ANTONIA (row 78): clear L, M, X flags
$ws.Range("D78").Value = ""
$ws.Range("E78").Value = ""
$ws.Range("F78").Value = ""

# MIGUEL (row 80): clear J flag
$ws.Range("G80").Value = ""

# LAURA (row 82): clear M flag
$ws.Range("E82").Value = ""

# Restore the default view for the sheet (scroll position / selection)
$ws.Range("I90").Select()
